$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1607',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1607',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1608',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1608',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1609',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1609',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1610',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1610',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1611',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1611',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1612',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1612',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1701',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1701',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1702',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1702',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1703',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1703',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1704',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1704',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1705',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1705',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1706',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1706',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1707',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1707',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1708',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1708',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1709',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1709',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1710',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1710',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1711',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1711',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1712',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1712',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1801',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1801',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1802',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1802',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1803',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1803',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1804',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1804',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1805',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1805',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1806',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1806',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1807',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1807',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1808',27580,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1808',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1809',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1809',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1810',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1810',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1811',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1811',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1812',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1812',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1901',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1901',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1902',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1902',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1903',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1903',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1904',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1904',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1905',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1905',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1906',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1906',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1907',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1907',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1908',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1908',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1909',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1909',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1910',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1910',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1911',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1911',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','1912',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','1912',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2001',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2001',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2002',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2002',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2003',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2003',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2004',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2004',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2005',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2005',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2006',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2006',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2007',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2007',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2008',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2008',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2009',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2009',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2010',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2010',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2011',31249,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2011',27578,689455),
    @('CC','1050971273','YAIR IVAN MEZA MARRIAGA','2012',21874,781242),
    @('CC','1003191558','CESAR LUIS PACHECO RAMIREZ','2012',19305,689455)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 16 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}